$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103001
$ws.Cells.Item($row, 10).Value = "Cereza"
$ws.Cells.Item($row, 11).Value = "Lapins"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 400
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 15500
$ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 1550
$ws.Cells.Item($row, 20).Value = 10
